# Update Ticket Sales (Q) and Embarking (R) figures for station rows
# as part of "updates to time module".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;   Q = 36; R = 31 },
    @{ Row = 10;  Q = 38; R = 35 },
    @{ Row = 17;  Q = 29; R = 25 },
    @{ Row = 23;  Q = 60; R = 31 },
    @{ Row = 32;  Q = 21; R = 21 },
    @{ Row = 40;  Q = 12; R = 10 },
    @{ Row = 49;  Q = 50; R = 1  },
    @{ Row = 58;  Q = 4;  R = 4  },
    @{ Row = 66;  Q = 57; R = 19 },
    @{ Row = 74;  Q = 48; R = 15 },
    @{ Row = 78;  Q = 3;  R = 2  },
    @{ Row = 89;  Q = 68; R = 62 },
    @{ Row = 97;  Q = 82; R = 26 },
    @{ Row = 106; Q = 16; R = 13 },
    @{ Row = 115; Q = 7;  R = 6  },
    @{ Row = 124; Q = 12; R = 9  },
    @{ Row = 133; Q = 61; R = 46 },
    @{ Row = 142; Q = 5;  R = 1  }
)

foreach ($u in $updates) {
    $ws.Range("Q$($u.Row)").Value = $u.Q
    $ws.Range("R$($u.Row)").Value = $u.R
}

$wb.Save()
